$wb = $excel.ActiveWorkbook

# --- ALC sheet updates ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 427911.12
$ws.Range("I28").Value = 653820.5
$ws.Range("J28").Value = 1193.3334
$ws.Range("K28").Value = 653820.5
$ws.Range("L28").Value = 1193.3334
$ws.Range("M28").Value = -653335.5
$ws.Range("N28").Value = -2163.3334
$ws.Range("H96").Value = 484.16666
$ws.Range("I96").Value = 432
$ws.Range("J96").Value = 588.5
$ws.Range("K96").Value = 1296
$ws.Range("L96").Value = 1765.5
$ws.Range("M96").Value = 77
$ws.Range("N96").Value = -4511.5
$ws.Range("H100").Value = 7576367
$ws.Range("I100").Value = 11111642
$ws.Range("J100").Value = 778.5714
$ws.Range("K100").Value = 11111642
$ws.Range("L100").Value = 778.5714
$ws.Range("M100").Value = -11111101
$ws.Range("N100").Value = -1860.5714
$ws.Range("H112").Value = 38961870
$ws.Range("I112").Value = 300
$ws.Range("J112").Value = 45455464
$ws.Range("K112").Value = 900
$ws.Range("L112").Value = 136366392
$ws.Range("M112").Value = 208
$ws.Range("N112").Value = -136368608
$ws.Range("H113").Value = 157500.72
$ws.Range("I113").Value = 178750.83
$ws.Range("J113").Value = 30000
$ws.Range("K113").Value = 178750.83
$ws.Range("L113").Value = 30000
$ws.Range("M113").Value = -175496.83
$ws.Range("N113").Value = -36508
$ws.Range("H121").Value = 1000
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 1000
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 3000
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -6494
$ws.Range("H133").Value = 11336
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 11336
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 11336
$ws.Range("N133").Value = -21456
$ws.Range("H138").Value = 5146056.5
$ws.Range("I138").Value = 1425111
$ws.Range("J138").Value = 7249200
$ws.Range("K138").Value = 4275333
$ws.Range("L138").Value = 21747600
$ws.Range("M138").Value = -4270193
$ws.Range("N138").Value = -21757880
$ws.Range("H141").Value = 1967.7805
$ws.Range("I141").Value = 1435.6061
$ws.Range("J141").Value = 4163
$ws.Range("K141").Value = 4306.8183
$ws.Range("L141").Value = 12489
$ws.Range("M141").Value = 873.1817000000001
$ws.Range("N141").Value = -22849

# --- ARM sheet updates ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14312.85
$ws.Range("I32").Value = 1604.3611
$ws.Range("J32").Value = 128689.25
$ws.Range("K32").Value = 1604.3611
$ws.Range("L32").Value = 128689.25
$ws.Range("M32").Value = -1317.3611
$ws.Range("N32").Value = -129263.25
$ws.Range("H61").Value = 3350.1052
$ws.Range("I61").Value = 2140.7273
$ws.Range("J61").Value = 5013
$ws.Range("K61").Value = 2140.7273
$ws.Range("L61").Value = 5013
$ws.Range("M61").Value = -1928.7273
$ws.Range("N61").Value = -5437
$ws.Range("H110").Value = 542.6
$ws.Range("I110").Value = 503.25
$ws.Range("K110").Value = 503.25
$ws.Range("M110").Value = 1541.75
$ws.Range("H132").Value = 2712.9375
$ws.Range("I132").Value = 2389.4119
$ws.Range("J132").Value = 3498.6428
$ws.Range("K132").Value = 7168.2357
$ws.Range("L132").Value = 10495.9284
$ws.Range("M132").Value = -4638.2357
$ws.Range("N132").Value = -15555.9284
$ws.Range("H133").Value = 38885.9
$ws.Range("J133").Value = 38885.9
$ws.Range("L133").Value = 38885.9
$ws.Range("N133").Value = -43945.9
$ws.Range("H136").Value = 3350.1052
$ws.Range("I136").Value = 2140.7273
$ws.Range("J136").Value = 5013
$ws.Range("K136").Value = 6422.1819
$ws.Range("L136").Value = 15039
$ws.Range("M136").Value = -3872.1819
$ws.Range("N136").Value = -20139
$ws.Range("H139").Value = 38527.6
$ws.Range("J139").Value = 38527.6
$ws.Range("L139").Value = 38527.6
$ws.Range("N139").Value = -48807.6

# --- BSM sheet updates ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 46625
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 46625
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 46625
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -48319
$ws.Range("H102").Value = 8800
$ws.Range("I102").Value = 8800
$ws.Range("K102").Value = 8800
$ws.Range("M102").Value = -5555
$ws.Range("H105").Value = 3583.3333
$ws.Range("I105").Value = 3375
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 3375
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -1628
$ws.Range("N105").Value = -7494
$ws.Range("H133").Value = 37776.25
$ws.Range("J133").Value = 37776.25
$ws.Range("L133").Value = 37776.25
$ws.Range("N133").Value = -47896.25

# --- CRP sheet updates ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("N44").ClearContents()

# --- CUL sheet updates ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 361.66666
$ws.Range("I50").Value = 292.5
$ws.Range("J50").Value = 500
$ws.Range("K50").Value = 877.5
$ws.Range("L50").Value = 1500
$ws.Range("M50").Value = -396.5
$ws.Range("N50").Value = -2462
$ws.Range("H53").Value = 361.66666
$ws.Range("I53").Value = 292.5
$ws.Range("J53").Value = 500
$ws.Range("K53").Value = 877.5
$ws.Range("L53").Value = 1500
$ws.Range("M53").Value = -396.5
$ws.Range("N53").Value = -2462
$ws.Range("H87").Value = 20333.25
$ws.Range("I87").Value = 8333
$ws.Range("J87").Value = 24333.334
$ws.Range("K87").Value = 24999
$ws.Range("L87").Value = 73000.00199999999
$ws.Range("M87").Value = -23751
$ws.Range("N87").Value = -75496.00199999999
$ws.Range("H90").Value = 20333.25
$ws.Range("I90").Value = 8333
$ws.Range("J90").Value = 24333.334
$ws.Range("K90").Value = 74997
$ws.Range("L90").Value = 219000.006
$ws.Range("M90").Value = -68757
$ws.Range("N90").Value = -231480.006

# --- GSM sheet updates ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7039.421
$ws.Range("I70").Value = 7364.2144
$ws.Range("K70").Value = 7364.2144
$ws.Range("M70").Value = -7094.2144
$ws.Range("H73").Value = 7039.421
$ws.Range("I73").Value = 7364.2144
$ws.Range("K73").Value = 7364.2144
$ws.Range("M73").Value = -6428.2144
$ws.Range("H137").Value = 45000
$ws.Range("J137").Value = 45000
$ws.Range("L137").Value = 45000
$ws.Range("N137").Value = -55200
$ws.Range("H138").Value = 49431.11
$ws.Range("J138").Value = 49431.11
$ws.Range("L138").Value = 49431.11
$ws.Range("N138").Value = -59711.11
$ws.Range("H139").Value = 35313
$ws.Range("J139").Value = 35313
$ws.Range("L139").Value = 35313
$ws.Range("N139").Value = -45593

# --- LTW sheet updates ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 952.5714
$ws.Range("I46").Value = 885.7143
$ws.Range("J46").Value = 1086.2858
$ws.Range("K46").Value = 885.7143
$ws.Range("L46").Value = 1086.2858
$ws.Range("M46").Value = -697.7143
$ws.Range("N46").Value = -1462.2858
$ws.Range("H132").Value = 4657.3057
$ws.Range("I132").Value = 4027.4583
$ws.Range("J132").Value = 5917
$ws.Range("K132").Value = 12082.3749
$ws.Range("L132").Value = 17751
$ws.Range("M132").Value = -9552.374899999999
$ws.Range("N132").Value = -22811
$ws.Range("H136").Value = 3528.7827
$ws.Range("I136").Value = 1720.2759
$ws.Range("J136").Value = 6613.8823
$ws.Range("K136").Value = 5160.8277
$ws.Range("L136").Value = 19841.6469
$ws.Range("M136").Value = -2610.8277
$ws.Range("N136").Value = -24941.6469

# --- WVR sheet updates ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 250000510
$ws.Range("I96").Value = 500000400
$ws.Range("J96").Value = 650
$ws.Range("K96").Value = 500000400
$ws.Range("L96").Value = 650
$ws.Range("M96").Value = -499999027
$ws.Range("N96").Value = -3396
$ws.Range("H132").Value = 13515718
$ws.Range("I132").Value = 29414188
$ws.Range("K132").Value = 88242564
$ws.Range("M132").Value = -88240034
$ws.Range("H138").Value = 70000
$ws.Range("J138").Value = 70000
$ws.Range("L138").Value = 70000
$ws.Range("N138").Value = -80280
